$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "timestamp" column (O) for every data row (2-63): the scrape
# that produced this workbook re-ran later the same day, so every row's
# captured timestamp moves from 07:10:23 to 21:00:44.
$ws.Range("O2:O63").Value = "2022-09-16 21:00:44"
